$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values.
# Pure-numeric looking price strings are prefixed with a leading apostrophe so
# Excel stores them as literal text (matching the original inlineStr content)
# instead of auto-converting them to numeric values.

$ws.Range("D2").Value = '27.018.46'
$ws.Range("E2").Value = '  +5.58%  '
$ws.Range("D3").Value = '1.879.15'
$ws.Range("E3").Value = '  +4.28%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("D5").Value = '''281.17'
$ws.Range("E5").Value = '  +2.72%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '''0.5263'
$ws.Range("E7").Value = '  +4.93%  '
$ws.Range("D8").Value = '''0.3534'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").Value = '''0.07036'
$ws.Range("E9").Value = '  +6.67%  '
$ws.Range("D10").Value = '''20.30'
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("D11").Value = '''0.8153'
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '''0.07789'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = '1.881.69'
$ws.Range("E13").Value = '  +4.46%  '
$ws.Range("D14").Value = '''5.221'
$ws.Range("E14").Value = '  +3.60%  '
$ws.Range("D15").Value = '''90.45'
$ws.Range("D16").Value = '''1.0000'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").Value = '''14.59'
$ws.Range("E17").Value = '  +5.36%  '
$ws.Range("D18").Value = '''0.000008169'
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("D19").Value = '''0.9998'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '27.047.04'
$ws.Range("E20").Value = '  +5.44%  '
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("D22").Value = '''10.19'
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("D23").Value = '''6.234'
$ws.Range("E23").Value = '  +3.27%  '
$ws.Range("D24").Value = '''2.390'
$ws.Range("E24").Value = '  +13.94%  '
$ws.Range("D25").Value = '''146.71'
$ws.Range("E25").Value = '  +3.67%  '
$ws.Range("D26").Value = '''17.59'
$ws.Range("E26").Value = '  +4.24%  '
$ws.Range("D27").Value = '''1.677'
$ws.Range("E27").Value = '  +1.47%  '
$ws.Range("D28").Value = '''113.43'
$ws.Range("E28").Value = '  +5.13%  '
$ws.Range("D29").Value = '''4.376'
$ws.Range("E29").Value = '  +1.86%  '
$ws.Range("D30").Value = '''4.377'
$ws.Range("E30").Value = '  +4.96%  '
$ws.Range("D31").Value = '''0.08886'
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("D32").Value = '''0.04897'
$ws.Range("D33").Value = '''1.172'
$ws.Range("E33").Value = '  +4.33%  '
$ws.Range("D34").Value = '''0.7414'
$ws.Range("E34").Value = '  +3.32%  '
$ws.Range("D35").Value = '''2.877'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = '''3.303'
$ws.Range("E36").Value = '  +9.39%  '
$ws.Range("D37").Value = '''2.410'
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("D38").Value = '''0.5300'
$ws.Range("E38").Value = '  +3.83%  '
$ws.Range("D39").Value = '''0.01881'
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("D40").Value = '''0.9813'
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("E41").Value = '  +2.86%  '
$ws.Range("D42").Value = '''6.316'
$ws.Range("E42").Value = '  +2.73%  '
$ws.Range("D43").Value = '''8.170'
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '''0.9998'
$ws.Range("D45").Value = '''0.4600'
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("D46").Value = '''0.1365'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '''9.457'
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("D48").Value = '''36.71'
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("D49").Value = '''1.519'
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").Value = '''61.81'
$ws.Range("E51").Value = '  +4.65%  '
